# Apply updated market-price derived values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 976.3333
$ws.Range("I18").Value = 976.3333
$ws.Range("K18").Value = 976.3333
$ws.Range("M18").Value = -692.3333
$ws.Range("H28").Value = 1771.3077
$ws.Range("I28").Value = 1473.3158
$ws.Range("K28").Value = 1473.3158
$ws.Range("M28").Value = -988.3158000000001
$ws.Range("H32").Value = 8355.223
$ws.Range("I32").Value = 1249.8334
$ws.Range("J32").Value = 22566
$ws.Range("K32").Value = 1249.8334
$ws.Range("L32").Value = 22566
$ws.Range("M32").Value = -923.8334
$ws.Range("N32").Value = -23218
$ws.Range("H74").Value = 12594.579
$ws.Range("I74").Value = 12594.579
$ws.Range("K74").Value = 12594.579
$ws.Range("M74").Value = -11658.579
$ws.Range("H77").Value = 12594.579
$ws.Range("I77").Value = 12594.579
$ws.Range("K77").Value = 62972.895
$ws.Range("M77").Value = -58292.895
$ws.Range("H86").Value = 1889.2727
$ws.Range("I86").Value = 2161.625
$ws.Range("J86").Value = 1163
$ws.Range("K86").Value = 2161.625
$ws.Range("L86").Value = 1163
$ws.Range("M86").Value = -1038.625
$ws.Range("N86").Value = -3409
$ws.Range("H89").Value = 1889.2727
$ws.Range("I89").Value = 2161.625
$ws.Range("J89").Value = 1163
$ws.Range("K89").Value = 10808.125
$ws.Range("L89").Value = 5815
$ws.Range("M89").Value = -5192.125
$ws.Range("N89").Value = -17047
$ws.Range("H106").Value = 9595.362999999999
$ws.Range("I106").Value = 9595.362999999999
$ws.Range("K106").Value = 9595.362999999999
$ws.Range("M106").Value = -8964.362999999999
$ws.Range("H132").Value = 3657.932
$ws.Range("I132").Value = 3565.5476
$ws.Range("J132").Value = 5598
$ws.Range("K132").Value = 10696.6428
$ws.Range("L132").Value = 16794
$ws.Range("M132").Value = -8166.6428
$ws.Range("N132").Value = -21854
$ws.Range("H138").Value = 4186.0894
$ws.Range("J138").Value = 4946.1113
$ws.Range("L138").Value = 14838.3339
$ws.Range("N138").Value = -25118.3339

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7028.853
$ws.Range("I32").Value = 2362.3333
$ws.Range("K32").Value = 2362.3333
$ws.Range("M32").Value = -2075.3333
$ws.Range("H45").Value = 3164.3333
$ws.Range("I45").Value = 2970.5334
$ws.Range("K45").Value = 2970.5334
$ws.Range("M45").Value = -2593.5334
$ws.Range("H74").Value = 5288.25
$ws.Range("I74").Value = 2662.2
$ws.Range("K74").Value = 2662.2
$ws.Range("M74").Value = -1788.2
$ws.Range("H77").Value = 5288.25
$ws.Range("I77").Value = 2662.2
$ws.Range("K77").Value = 13311
$ws.Range("M77").Value = -8943
$ws.Range("H97").Value = 856.5333000000001
$ws.Range("I97").Value = 856.5333000000001
$ws.Range("K97").Value = 856.5333000000001
$ws.Range("M97").Value = -360.5333000000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3204
$ws.Range("I94").Value = 1422.6428
$ws.Range("K94").Value = 1422.6428
$ws.Range("M94").Value = -971.6428000000001
$ws.Range("H134").Value = 2278.068
$ws.Range("I134").Value = 2278.068
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6834.204000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4299.204000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 26581.084
$ws.Range("I16").Value = 16270.363
$ws.Range("K16").Value = 16270.363
$ws.Range("M16").Value = -15983.363
$ws.Range("H31").Value = 6602.7
$ws.Range("I31").Value = 6684.4443
$ws.Range("J31").Value = 6480.0835
$ws.Range("K31").Value = 6684.4443
$ws.Range("L31").Value = 6480.0835
$ws.Range("M31").Value = -6389.4443
$ws.Range("N31").Value = -7070.0835
$ws.Range("H34").Value = 6602.7
$ws.Range("I34").Value = 6684.4443
$ws.Range("J34").Value = 6480.0835
$ws.Range("K34").Value = 6684.4443
$ws.Range("L34").Value = 6480.0835
$ws.Range("M34").Value = -6482.4443
$ws.Range("N34").Value = -6884.0835
$ws.Range("H113").Value = 26581.084
$ws.Range("I113").Value = 16270.363
$ws.Range("K113").Value = 16270.363
$ws.Range("M113").Value = -14100.363
$ws.Range("H131").Value = 53327.332
$ws.Range("J131").Value = 55994.2
$ws.Range("L131").Value = 55994.2
$ws.Range("N131").Value = -66074.2
$ws.Range("H132").Value = 3325.5715
$ws.Range("I132").Value = 2995.75
$ws.Range("K132").Value = 8987.25
$ws.Range("M132").Value = -6457.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1280.4166
$ws.Range("J113").Value = 1523.3889
$ws.Range("L113").Value = 4570.1667
$ws.Range("N113").Value = -8910.1667

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 32971.6
$ws.Range("J15").Value = 32971.6
$ws.Range("L15").Value = 32971.6
$ws.Range("N15").Value = -33547.6
$ws.Range("H80").Value = 3951.2354
$ws.Range("I80").Value = 2715.9092
$ws.Range("J80").Value = 6216
$ws.Range("K80").Value = 2715.9092
$ws.Range("L80").Value = 6216
$ws.Range("M80").Value = -1717.9092
$ws.Range("N80").Value = -8212
$ws.Range("H81").Value = 32971.6
$ws.Range("J81").Value = 32971.6
$ws.Range("L81").Value = 32971.6
$ws.Range("N81").Value = -34967.6
$ws.Range("H83").Value = 3951.2354
$ws.Range("I83").Value = 2715.9092
$ws.Range("J83").Value = 6216
$ws.Range("K83").Value = 13579.546
$ws.Range("L83").Value = 31080
$ws.Range("M83").Value = -8587.546
$ws.Range("N83").Value = -41064
$ws.Range("H84").Value = 32971.6
$ws.Range("J84").Value = 32971.6
$ws.Range("L84").Value = 98914.79999999999
$ws.Range("N84").Value = -108898.8

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1666.3334
$ws.Range("I16").Value = 1666.3334
$ws.Range("K16").Value = 1666.3334
$ws.Range("M16").Value = -1496.3334
$ws.Range("H22").Value = 6591.6665
$ws.Range("I22").Value = 4421.3335
$ws.Range("J22").Value = 10932.333
$ws.Range("K22").Value = 4421.3335
$ws.Range("L22").Value = 10932.333
$ws.Range("M22").Value = -4126.3335
$ws.Range("N22").Value = -11522.333
$ws.Range("H27").Value = 6591.6665
$ws.Range("I27").Value = 4421.3335
$ws.Range("J27").Value = 10932.333
$ws.Range("K27").Value = 4421.3335
$ws.Range("L27").Value = 10932.333
$ws.Range("M27").Value = -4314.3335
$ws.Range("N27").Value = -11146.333
$ws.Range("H55").Value = 1064.4667
$ws.Range("I55").Value = 972.3684
$ws.Range("J55").Value = 1223.5454
$ws.Range("K55").Value = 972.3684
$ws.Range("L55").Value = 1223.5454
$ws.Range("M55").Value = -799.3684
$ws.Range("N55").Value = -1569.5454
$ws.Range("H82").Value = 1986.4
$ws.Range("I82").Value = 2033
$ws.Range("K82").Value = 2033
$ws.Range("M82").Value = -1672
$ws.Range("H85").Value = 1986.4
$ws.Range("I85").Value = 2033
$ws.Range("K85").Value = 2033
$ws.Range("M85").Value = -785
$ws.Range("H132").Value = 6274.2666
$ws.Range("I132").Value = 3222.4285
$ws.Range("K132").Value = 9667.2855
$ws.Range("M132").Value = -7137.2855

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3621.3333
$ws.Range("I62").Value = 4099.5
$ws.Range("K62").Value = 4099.5
$ws.Range("M62").Value = -3475.5
$ws.Range("H65").Value = 3621.3333
$ws.Range("I65").Value = 4099.5
$ws.Range("K65").Value = 20497.5
$ws.Range("M65").Value = -17377.5
$ws.Range("H96").Value = 3898.8572
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("N96").Value = -7746
$ws.Range("H100").Value = 2526.5715
$ws.Range("I100").Value = 2114.3333
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 4228.6666
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -3687.6666
$ws.Range("N100").Value = -11082
$ws.Range("H107").Value = 6496.6206
$ws.Range("I107").Value = 6360.316
$ws.Range("K107").Value = 19080.948
$ws.Range("M107").Value = -17160.948
$ws.Range("H132").Value = 3968.9556
$ws.Range("I132").Value = 3927.75
$ws.Range("K132").Value = 11783.25
$ws.Range("M132").Value = -9253.25

# ---- Cell removals ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N134").ClearContents()
